# Generate Report for Handback
# -----------------------------------------------------------------------
# Populates the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns for the two localized-language report sheets
# (zh-cn, de-de) now that the de-de handback has come back in sync with
# en-US, and flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is shown (the per-row
# detail sheets as well as the Overview roll-up). It also widens a couple
# of columns that now hold longer file-name/timestamp content.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$sourceMdDisplay = "89e05de2-83ad-4e5e-a67d-1b6105ce6083.md"
$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7dbdaca14229b83a9d6df9f0584f1e7a28c00dd8/e2e/89e05de2-83ad-4e5e-a67d-1b6105ce6083.md"

$zhHandbackFile = "89e05de2-83ad-4e5e-a67d-1b6105ce6083.12379d27f82a64c2b63665bf50a6c8754ebed751.zh-cn.xlf"
$deHandbackFile = "89e05de2-83ad-4e5e-a67d-1b6105ce6083.12379d27f82a64c2b63665bf50a6c8754ebed751.de-de.xlf"

$zhHandbackDateTime = "2016-08-23 13:04:36"
$deHandbackDateTime = "2016-08-23 13:04:43"

# ---------------------------------------------------------------------
# Overview sheet: roll-up Status cells for both languages (cols E/F),
# plus widen those two columns to fit the new, longer status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (C) widens to fit the new text.
$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527

# Latest Target File (I): hyperlink to the source markdown, same as col A.
$zh.Hyperlinks.Add($zh.Range("I2"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null

# Latest Handback File (J) / Latest Handback DateTime (K).
$zh.Range("J2").Value = $zhHandbackFile
$zh.Range("J3").Value = $zhHandbackFile
$zh.Range("K2").Value = $zhHandbackDateTime
$zh.Range("K3").Value = $zhHandbackDateTime

# Latest Target File / Latest Handback File columns widen.
$zh.Columns.Item(9).ColumnWidth = 39.1666666666667
$zh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column (C) widens to fit the new text.
$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew
$de.Columns.Item(3).ColumnWidth = 29.9777047293527

# Latest Target File (I): hyperlink to the source markdown, same as col A.
$de.Hyperlinks.Add($de.Range("I2"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null

# Latest Handback File (J) / Latest Handback DateTime (K) - de-de is now
# in sync, so it gets its own xlf + the newest handback timestamp.
$de.Range("J2").Value = $deHandbackFile
$de.Range("J3").Value = $deHandbackFile
$de.Range("K2").Value = $deHandbackDateTime
$de.Range("K3").Value = $deHandbackDateTime

# Latest Target File / Latest Handback File columns widen.
$de.Columns.Item(9).ColumnWidth = 39.1666666666667
$de.Columns.Item(10).ColumnWidth = 39.1666666666667
